$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row to reflect the new column names (handles multiple primary keys)
$ws.Range("A1").Value = "TradeID"
$ws.Range("B1").Value = "Risk"
$ws.Range("C1").Value = "Curve"
$ws.Range("D1").Value = "Type"

# Move the active selection to B1
$ws.Range("B1").Select()
